$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "corresponding threshold"
$ws.Range("B9").Value = "corresponding accuracy"
$ws.Range("B10").Value = "corresponding True Positive Rate"
$ws.Range("B11").Value = "corresponding False Positive Rate"
$ws.Range("B12").Value = "corresponding precision"
